$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-0.979***"
$ws.Range("C2").Value = "0.03***"
$ws.Range("D2").Value = "-0.028***"

$ws.Range("B3").Value = "-2.187**"
$ws.Range("C3").Value = "-0.195***"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.075"
$ws.Range("D3").Style = "Normal"

$ws.Range("B4").Value = "11.552***"
$ws.Range("C4").Value = "-1.066***"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.118"
$ws.Range("D4").Style = "Normal"

$ws.Range("A5:D6").Delete()
